$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.720.84"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "2.091.61"
$ws.Range("E3").Value = "  +2.68%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "228.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.612"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.88%  "

$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.384"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0832"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.69%  "

$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("D12").Value = "2.407.06"
$ws.Range("E12").Value = "  +2.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "22.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.89%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.797"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.45"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.04%  "

$ws.Range("D17").Value = "2.094.16"
$ws.Range("E17").Value = "  +2.74%  "

$ws.Range("D18").Value = "38.790.01"
$ws.Range("E18").Value = "  +2.87%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "71.91"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.26%  "

$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.29%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "226.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.05%  "

$ws.Range("E23").Value = "  -0.39%  "

$ws.Range("E24").Value = "  +1.31%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.33"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.62%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.40"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.31%  "

$ws.Range("E28").Value = "  +5.64%  "

$ws.Range("E29").Value = "  +9.27%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.15"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.90%  "

$ws.Range("E31").Value = "  +0.77%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.34"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.73"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +5.65%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.48"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0613"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.27%  "

$ws.Range("E36").Value = "  +2.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.50"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.57%  "

$ws.Range("E39").Value = "  +0.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.33"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "101.19"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.02%  "

$ws.Range("D42").Value = "1.538.89"
$ws.Range("E42").Value = "  +0.76%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0221"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.97%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0927"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.34%  "

$ws.Range("E45").Value = "  -0.94%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.66"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +9.48%  "

$ws.Range("E47").Value = "  +0.98%  "

$ws.Range("E48").Value = "  -4.33%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.83%  "

$ws.Range("E50").Value = "  +1.18%  "

$ws.Range("D51").Value = "2.293.32"
$ws.Range("E51").Value = "  +2.97%  "
